$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old header row (row 24), pushing the
# header + data table down by two rows (24->26 .. 35->37).
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(24).Insert()

# New text row 23: a note about going back through Git history.
$ws.Range("A23").Value = "Going back in Git commits, to check the files I forgot to check for peak double-counting"

# The "double count check passed?" column (L) now gets a new answer for
# the first five data rows (originally rows 25-29, now 27-31), noting
# that the check was done by looking back at Git history.
$note = "yes looking back at Git history of notebook"
$ws.Range("L27").Value = $note
$ws.Range("L28").Value = $note
$ws.Range("L29").Value = $note
$ws.Range("L30").Value = $note
$ws.Range("L31").Value = $note

# Column L needs to be widened to fit the new, much longer text.
$ws.Columns.Item(12).ColumnWidth = 74.08984375

# Leave the selection on the newly-noted row, matching the author's
# last on-screen interaction.
$ws.Rows.Item(25).Select()
